$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings (e.g. "1.004", "24.148.38")
# are preserved exactly as text, matching the inline string cells in the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.148.38'
$ws.Range("E2").Value = '  -3.21%  '
$ws.Range("D3").Value = '1.645.46'
$ws.Range("E3").Value = '  -3.29%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("D5").Value = '308.44'
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.3890'
$ws.Range("E7").Value = '  -1.96%  '
$ws.Range("D8").Value = '0.3875'
$ws.Range("E8").Value = '  -3.61%  '
$ws.Range("D9").Value = '1.005'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '1.366'
$ws.Range("E10").Value = '  -6.80%  '
$ws.Range("E11").Value = '  -6.95%  '
$ws.Range("D12").Value = '0.08481'
$ws.Range("E12").Value = '  -3.48%  '
$ws.Range("D13").Value = '24.21'
$ws.Range("E13").Value = '  -6.46%  '
$ws.Range("D14").Value = '7.168'
$ws.Range("E14").Value = '  -3.75%  '
$ws.Range("D15").Value = '0.00001290'
$ws.Range("E15").Value = '  -4.35%  '
$ws.Range("D16").Value = '7.528'
$ws.Range("E16").Value = '  -5.31%  '
$ws.Range("D17").Value = '1.650.05'
$ws.Range("E17").Value = '  -3.70%  '
$ws.Range("D18").Value = '94.91'
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("E19").Value = '  -3.22%  '
$ws.Range("E20").Value = '  +3.14%  '
$ws.Range("D21").Value = '6.975'
$ws.Range("E21").Value = '  -4.99%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '13.84'
$ws.Range("E23").Value = '  -3.51%  '
$ws.Range("D24").Value = '24.158.96'
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("D26").Value = '2.740'
$ws.Range("E26").Value = '  -6.64%  '
$ws.Range("D27").Value = '22.58'
$ws.Range("E27").Value = '  -4.69%  '
$ws.Range("D28").Value = '8.865'
$ws.Range("E28").Value = '  +6.79%  '
$ws.Range("D29").Value = '158.07'
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").Value = '142.48'
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("D31").Value = '5.409'
$ws.Range("E31").Value = '  -12.60%  '
$ws.Range("D32").Value = '2.449'
$ws.Range("E32").Value = '  -6.42%  '
$ws.Range("D33").Value = '1.830.97'
$ws.Range("E33").Value = '  -3.63%  '
$ws.Range("D34").Value = '7.093'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("D35").Value = '0.08087'
$ws.Range("E35").Value = '  -5.20%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02967'
$ws.Range("E36").Value = '  -5.16%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.9907'
$ws.Range("E37").Value = '  -4.42%  '
$ws.Range("D38").Value = '0.2710'
$ws.Range("E38").Value = '  -4.92%  '
$ws.Range("D39").Value = '0.09305'
$ws.Range("E39").Value = '  -2.38%  '
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").Value = '10.09'
$ws.Range("E41").Value = '  -7.14%  '
$ws.Range("D42").Value = '0.7648'
$ws.Range("E42").Value = '  -6.76%  '
$ws.Range("D43").Value = '13.12'
$ws.Range("E43").Value = '  -5.73%  '
$ws.Range("D44").Value = '16.22'
$ws.Range("E44").Value = '  -5.51%  '
$ws.Range("D45").Value = '2.499'
$ws.Range("E45").Value = '  -6.46%  '
$ws.Range("D46").Value = '0.6897'
$ws.Range("E46").Value = '  -6.37%  '
$ws.Range("D47").Value = '4.093'
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("D48").Value = '1.003'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").Value = '0.08420'
$ws.Range("E49").Value = '  -3.78%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '134.34'
$ws.Range("E50").Value = '  -3.28%  '
$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").Value = '1.265'
$ws.Range("E51").Value = '  -9.40%  '